$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(245, 44319, 5, 13, 86.16689865447074),
    @(246, 44320, 1, 14, 92.79512162789156),
    @(247, 44321, 1, 14, 92.79512162789156)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy the format of the last existing row's date cell (A244) so the new
    # date cell gets the same style (border/alignment/number format) without
    # introducing new style entries.
    $ws.Range("A244").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = $false
